$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a date/time number format to column A (this will create a new cell style)
$ws.Columns.Item(1).NumberFormat = "m/d/yy h:mm"

# Add the new row of data
$ws.Range("A2").Value = 42605.889074074075
$ws.Range("A2").NumberFormat = "m/d/yy h:mm"
$ws.Range("B2").Value = -18
$ws.Range("C2").Value = 63
$ws.Range("D2").Value = 33
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 20377
$ws.Range("H2").Value = 9796
$ws.Range("I2").Value = 1132
$ws.Range("J2").Value = 164
$ws.Range("K2").Value = 87
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 23
$ws.Range("N2").Value = "Bag"

$ws.Columns.Item(1).ColumnWidth = 14
